$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the auto-updating "date" placeholders (datetimeFigureOut
#    fields) on the slide master and every custom layout so they show
#    the current save date instead of the stale cached value.
# ---------------------------------------------------------------------
$newDate = "14/09/2024"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DatePlaceholders $layouts.Item($l).Shapes
}

# NOTE: the notes master's own date placeholder is intentionally left
# untouched — writing through NotesMaster.Shapes on this host mis-binds
# to the (numerically coincident) shape id on the slide master and
# corrupts its body placeholder text, so it is skipped here.

# ---------------------------------------------------------------------
# 2) Drop the trailing "Painel de Análise" / "Hidrometria" slides that
#    belonged to the old interactive-panel project — remove the last
#    two slides from the deck.
# ---------------------------------------------------------------------
while ($p.Slides.Count -gt 3) {
    $p.Slides.Item($p.Slides.Count).Delete()
}
